# A6-1.docx — rewrite the assignment body:
#   * shorten the bold title run from "A6-1: Unit testing" to "A6-1: "
#   * drop all the screenshot paragraphs ("Calculator class:", the two
#     Calculator screenshots, "UnitTest1 class:", the three UnitTest1
#     screenshots, "Test Explorer:" and its screenshot)
#   * replace the old "Calculator class:" heading paragraph with a new,
#     non-bold body paragraph describing the Task.WhenAll exercise
#
# The trailing empty paragraph (and the sectPr after it) is left alone.

$d = $word.ActiveDocument

$wNs   = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'
$w14Ns = 'http://schemas.microsoft.com/office/word/2010/wordml'

# ---------------------------------------------------------------------
# 1) Remove paragraphs 3..10 (1-based): the two Calculator screenshots,
#    "UnitTest1 class:", its three screenshots, "Test Explorer:" and its
#    screenshot. Paragraph.Range.Delete() also consumes the paragraph
#    mark, so repeatedly deleting whatever is now paragraph 3 removes
#    them one at a time without needing to recompute an end index.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 8; $i++) {
    $victim = $d.Paragraphs.Item(3)
    [void]$victim.Range.Delete()
}

# ---------------------------------------------------------------------
# 2) Shorten the title paragraph's run text, keeping its existing
#    paragraph/run identity (paraId/textId/rsid*) and bold/28pt
#    formatting intact.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titleXml = "<w:p xmlns:w='$wNs' xmlns:w14='$w14Ns' w14:paraId='6C132525' w14:textId='19AC0BE4' w:rsidR='00F42B13' w:rsidRDefault='006E6203'>" +
            "<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr></w:pPr>" +
            "<w:r w:rsidRPr='006E6203'><w:rPr><w:b/><w:bCs/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr>" +
            "<w:t xml:space='preserve'>A6-1: </w:t></w:r></w:p>"
[void]$titlePara.Range.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 3) Replace the old bold "Calculator class:" heading paragraph with the
#    new 12pt (sz 24) body paragraph describing the exercise. It is
#    built from three runs so "Task.WhenAll" can be wrapped in the
#    spell-check proofErr markers, matching the target markup.
# ---------------------------------------------------------------------
$bodyPara = $d.Paragraphs.Item(2)
$bodyXml = "<w:p xmlns:w='$wNs'>" +
           "<w:pPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr>" +
           "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr>" +
           "<w:t xml:space='preserve'>Modify the example of Fig. 23.3 of your C# book (pg. 953) to process the results of the tasks.  Use an array or list on the Task produced from the </w:t></w:r>" +
           "<w:proofErr w:type='spellStart'/>" +
           "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>Task.WhenAll</w:t></w:r>" +
           "<w:proofErr w:type='spellEnd'/>" +
           "<w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> method.</w:t></w:r>" +
           "</w:p>"
[void]$bodyPara.Range.InsertXML($bodyXml)
